# "Added new and updated report." - append the new LeetCode entry (#10,
# "house robber 3") to the Report sheet, then leave the Report sheet as
# the active tab with the selection parked on the next empty row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Report")

$ws.Range("A13").Value = 10
$ws.Range("C13").Value = "house robber 3"
$ws.Range("D13").Value = "Medium"
$ws.Range("E13").Value = "d&c"
$ws.Range("F13").Value = "Medium"
$ws.Range("G13").Value = "dp"
$ws.Range("H13").Value = 120

$ws.Activate()
$ws.Range("I14").Select()
